$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.803.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("E2").Style = "Normal"

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.149.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E3").Style = "Normal"

# Row 4: TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("E4").Style = "Normal"

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E5").Style = "Normal"

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.12%  "
$ws.Range("E6").Style = "Normal"

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E7").Style = "Normal"

# Row 8: LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.146.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E8").Style = "Normal"

# Row 9: XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("E9").Style = "Normal"

# Row 10: Dogecoin
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("E10").Style = "Normal"

# Row 11: Toncoin
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.78%  "
$ws.Range("E11").Style = "Normal"

# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("E12").Style = "Normal"

# Row 13: Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.78%  "
$ws.Range("E13").Style = "Normal"

# Row 14: ShibaInu
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000248"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("E14").Style = "Normal"

# Row 15: WrappedEther
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.504.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +11.84%  "
$ws.Range("E15").Style = "Normal"

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.665.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E16").Style = "Normal"

# Row 17: TRON
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.121"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("E17").Style = "Normal"

# Row 18: Polkadot
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("E18").Style = "Normal"

# Row 19: WrappedBTC
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.373.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("E19").Style = "Normal"

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("E20").Style = "Normal"

# Row 21: Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.57%  "
$ws.Range("E21").Style = "Normal"

# Row 22: Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.757"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.52%  "
$ws.Range("E22").Style = "Normal"

# Row 23: Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.51%  "
$ws.Range("E23").Style = "Normal"

# Row 24: InternetComputer(DFINITY)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.25%  "
$ws.Range("E24").Style = "Normal"

# Row 25: Fetch.AI
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.10%  "
$ws.Range("E25").Style = "Normal"

# Row 26: Litecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("E26").Style = "Normal"

# Row 27: Dai
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E27").Style = "Normal"

# Row 28: RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.20%  "
$ws.Range("E28").Style = "Normal"

# Row 29: PancakeSwap
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("E29").Style = "Normal"

# Row 30: NEARProtocol
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.52%  "
$ws.Range("E30").Style = "Normal"

# Row 31: FirstDigitalUSD
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("E31").Style = "Normal"

# Row 32: ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("E32").Style = "Normal"

# Row 33: Hedera
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.58%  "
$ws.Range("E33").Style = "Normal"

# Row 34: EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.25%  "
$ws.Range("E34").Style = "Normal"

# Row 35: PEPE
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0877"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.90%  "
$ws.Range("E35").Style = "Normal"

# Row 36: Mantle
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.39%  "
$ws.Range("E36").Style = "Normal"

# Row 37: dogwifhat
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.75%  "
$ws.Range("E37").Style = "Normal"

# Row 38: Filecoin
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.25%  "
$ws.Range("E38").Style = "Normal"

# Row 39: Stacks
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("E39").Style = "Normal"

# Row 40: Bittensor
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "463.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.64%  "
$ws.Range("E40").Style = "Normal"

# Row 41: Cosmos
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.46%  "
$ws.Range("E41").Style = "Normal"

# Row 42: OKB
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "51.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E42").Style = "Normal"

# Row 43: TheGraph
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.299"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +8.08%  "
$ws.Range("E43").Style = "Normal"

# Row 44: VeChain
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0375"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("E44").Style = "Normal"

# Row 45: Maker
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.890.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("E45").Style = "Normal"

# Row 46: Kaspa
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.110"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.85%  "
$ws.Range("E46").Style = "Normal"

# Row 47: Arweave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.27%  "
$ws.Range("E47").Style = "Normal"

# Row 48: Monero
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.33%  "
$ws.Range("E48").Style = "Normal"

# Row 49: InjectiveProtocol
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.88%  "
$ws.Range("E49").Style = "Normal"

# Row 50: ThetaToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.79%  "
$ws.Range("E50").Style = "Normal"

# Row 51: USDe
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.00%  "
$ws.Range("E51").Style = "Normal"
